$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-by-row refresh of the "Price" (D) and "Volume(1h)" (E) columns of the
# cryptocurrency table, as produced by the scheduled GitHub Actions job.
# D is only listed when the price text actually changed; E always changes.
$updates = @(
    @{ Row = 2; D = "59.307.02"; E = "  +0.02%  " },
    @{ Row = 3; D = "2.643.95"; E = "  -0.40%  " },
    @{ Row = 4; D = $null; E = "  +0.18%  " },
    @{ Row = 5; D = "528.79"; E = "  +0.76%  " },
    @{ Row = 6; D = "145.56"; E = "  -0.59%  " },
    @{ Row = 7; D = $null; E = "  -0.07%  " },
    @{ Row = 8; D = "0.571"; E = "  -0.63%  " },
    @{ Row = 9; D = "6.67"; E = "  -3.72%  " },
    @{ Row = 10; D = $null; E = "  +1.04%  " },
    @{ Row = 11; D = $null; E = "  +0.25%  " },
    @{ Row = 12; D = $null; E = "  +0.57%  " },
    @{ Row = 13; D = "3.111.92"; E = "  -0.12%  " },
    @{ Row = 14; D = "59.314.03"; E = "  -0.01%  " },
    @{ Row = 15; D = "20.83"; E = "  -2.09%  " },
    @{ Row = 16; D = $null; E = "  -0.14%  " },
    @{ Row = 17; D = "2.635.22"; E = "  -0.47%  " },
    @{ Row = 18; D = "342.34"; E = "  -0.06%  " },
    @{ Row = 19; D = $null; E = "  +0.06%  " },
    @{ Row = 20; D = "10.62"; E = "  +1.97%  " },
    @{ Row = 21; D = "6.39"; E = "  +0.93%  " },
    @{ Row = 22; D = "0.999"; E = "  -0.19%  " },
    @{ Row = 23; D = "65.64"; E = "  +3.04%  " },
    @{ Row = 24; D = $null; E = "  +0.76%  " },
    @{ Row = 25; D = "0.167"; E = "  -0.52%  " },
    @{ Row = 26; D = "0.998"; E = "  -0.07%  " },
    @{ Row = 27; D = "7.20"; E = "  +0.95%  " },
    @{ Row = 28; D = "0.0₃0802"; E = "  -0.99%  " },
    @{ Row = 29; D = $null; E = "  -0.06%  " },
    @{ Row = 30; D = "6.39"; E = "  -4.80%  " },
    @{ Row = 31; D = $null; E = "  +0.80%  " },
    @{ Row = 32; D = "19.01"; E = "  +0.97%  " },
    @{ Row = 33; D = "150.49"; E = "  +0.84%  " },
    @{ Row = 34; D = "4.15"; E = "  -1.42%  " },
    @{ Row = 35; D = "1.19"; E = "  -1.08%  " },
    @{ Row = 36; D = "0.865"; E = "  -4.60%  " },
    @{ Row = 37; D = "0.861"; E = "  -3.20%  " },
    @{ Row = 38; D = "1.48"; E = "  -0.59%  " },
    @{ Row = 39; D = "36.52"; E = "  -0.78%  " },
    @{ Row = 40; D = $null; E = "  +1.23%  " },
    @{ Row = 41; D = $null; E = "  +0.14%  " },
    @{ Row = 42; D = $null; E = "  -0.15%  " },
    @{ Row = 43; D = "0.602"; E = "  -3.55%  " },
    @{ Row = 44; D = "270.09"; E = "  -2.10%  " },
    @{ Row = 45; D = "19.39"; E = "  -2.41%  " },
    @{ Row = 46; D = "10.71"; E = "  +1.76%  " },
    @{ Row = 47; D = "0.0535"; E = "  -0.93%  " },
    @{ Row = 48; D = "2.038.93"; E = "  -1.07%  " },
    @{ Row = 49; D = "4.76"; E = "  -2.46%  " },
    @{ Row = 50; D = $null; E = "  -0.73%  " },
    @{ Row = 51; D = "18.71"; E = "  -2.61%  " }
)

# Both columns hold plain text (e.g. "59.307.02", "0.571", "  +0.02%  ").
# Excel auto-detects plain numeric-looking strings (like "528.79") and would
# silently convert them to numbers (normalizing "7.20" -> 7.2 in the
# process), so text that would otherwise parse as a plain number is entered
# with a leading apostrophe to force text storage. Non-numeric text
# (multi-dot "59.307.02", space-padded percentages, etc.) is assigned as-is.
$numberPattern = '^[+-]?[0-9]*\.?[0-9]+$'

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($u.Row, 4)
        if ($u.D -match $numberPattern) {
            $dCell.Value = "'" + $u.D
        } else {
            $dCell.Value = $u.D
        }
    }

    $eCell = $ws.Cells.Item($u.Row, 5)
    if ($u.E -match $numberPattern) {
        $eCell.Value = "'" + $u.E
    } else {
        $eCell.Value = $u.E
    }
}
